# Applies the "Fri Sep 29 05:56:21 UTC 2023" cryptos list refresh described
# by the commit diff. Every cell in this sheet holds plain text (prices and
# percent-volume strings are formatted text, not numeric Excel values), so
# string literals are assigned directly. For new values that look like plain
# numbers (e.g. "215.20", "0.0618") a leading apostrophe forces Excel to keep
# them as text instead of silently re-interpreting them as numeric values;
# the style is then reset to "Normal" so no stray number-format is left on
# the cell (matching the original, unformatted text cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.980.63"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "1.655.87"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.Value = "'215.20"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E7").Value = "  -0.06%  "
$c = $ws.Range("D8")
$c.Value = "'0.251"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.56%  "
$c = $ws.Range("D9")
$c.Value = "'0.0618"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.89%  "
$c = $ws.Range("D10")
$c.Value = "'20.16"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.80%  "
$c = $ws.Range("D11")
$c.Value = "'0.0879"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "1.651.11"
$ws.Range("E13").Value = "  +2.60%  "
$c = $ws.Range("D14")
$c.Value = "'4.09"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "
$c = $ws.Range("D15")
$c.Value = "'0.524"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.88%  "
$c = $ws.Range("D16")
$c.Value = "'65.28"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "26.984.01"
$ws.Range("E17").Value = "  +2.10%  "
$c = $ws.Range("D18")
$c.Value = "'237.02"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "0.0₃0741"
$ws.Range("E19").Value = "  +2.17%  "
$c = $ws.Range("D20")
$c.Value = "'7.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("E21").Value = "  +0.00%  "
$c = $ws.Range("D22")
$c.Value = "'4.44"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("E23").Value = "  +2.94%  "
$c = $ws.Range("D24")
$c.Value = "'2.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.24%  "
$c = $ws.Range("D25")
$c.Value = "'145.33"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$c = $ws.Range("D26")
$c.Value = "'7.13"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  +0.89%  "
$c = $ws.Range("D28")
$c.Value = "'15.88"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("E29").Value = "  -0.02%  "
$c = $ws.Range("D30")
$c.Value = "'0.0498"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "1.554.32"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  +4.45%  "
$c = $ws.Range("D35")
$c.Value = "'1.62"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.99%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +3.89%  "
$c = $ws.Range("D38")
$c.Value = "'0.906"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D42")
$c.Value = "'66.65"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +8.71%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D43")
$c.Value = "'0.980"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.76%  "
$c = $ws.Range("D44")
$c.Value = "'2.24"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").Value = "1.797.03"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("E46").Value = "  +1.79%  "
$c = $ws.Range("D47")
$c.Value = "'90.24"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("E48").Value = "  +3.09%  "
$c = $ws.Range("D49")
$c.Value = "'0.0995"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.66%  "
$c = $ws.Range("D50")
$c.Value = "'0.0506"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  +2.26%  "
